# Update the dSF (column F) values as per repulled data.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$changes = @{
    2  = 6
    3  = 4
    4  = -1
    7  = 4
    8  = 5
    9  = -6
    10 = 5
    11 = -1
    12 = 3
    15 = 2
    16 = 5
    17 = -2
    18 = -1
    19 = 2
    21 = -2
    22 = 4
    23 = 3
    26 = -1
    27 = 3
    28 = -2
    29 = 1
    31 = 2
    32 = 8
    33 = 2
    34 = 3
    35 = 3
    36 = -3
    37 = 5
    39 = 1
}

foreach ($row in $changes.Keys) {
    $ws.Cells.Item($row, 6).Value = $changes[$row]
}
